$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row: insert two new columns (Microdol1, Microdol5) between
#     existing "SK2" column and shift it to column G, renaming it "SK2 (old data)" ---
$ws.Range("E1:F1").EntireColumn.Insert()

$ws.Range("E1").Value = "Microdol1"
$ws.Range("F1").Value = "Microdol5"
$ws.Range("G1").Value = "SK2 (old data)"

# --- Updated test data (October 2023 column test data) ---

# Row 2 - CaPct
$ws.Range("C2").Value = 39.6
$ws.Range("D2").Value = 39.6
$ws.Range("E2").Value = 21.5
$ws.Range("F2").Value = 21.5

# Row 3 - MgPct
$ws.Range("C3").Value = 0.4
$ws.Range("D3").Value = 0.4
$ws.Range("E3").Value = 12.8
$ws.Range("F3").Value = 12.8

# Row 4 - DryFac
$ws.Range("E4").Value = 0.7
$ws.Range("F4").Value = 0.7

# Row 5 - ColDepth
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 2

# Row 6 - IDph40
$ws.Range("C6").Value = 64
$ws.Range("D6").Value = 61.2
$ws.Range("E6").Value = 73.1
$ws.Range("F6").Value = 70.6

# Row 7 - IDph45
$ws.Range("C7").Value = 48.6
$ws.Range("D7").Value = 60.5
$ws.Range("E7").Value = 76.2
$ws.Range("F7").Value = 46.1

# Row 8 - IDph50
$ws.Range("C8").Value = 50.3
$ws.Range("D8").Value = 60.2
$ws.Range("E8").Value = 25.1
$ws.Range("F8").Value = 33.4

# Row 9 - IDph55
$ws.Range("C9").Value = 43
$ws.Range("D9").Value = 51.7
$ws.Range("E9").Value = 25.6
$ws.Range("F9").Value = 25.3

# Row 10 - IDph60
$ws.Range("C10").Value = 34.9
$ws.Range("D10").Value = 48.5
$ws.Range("E10").Value = 24.8
$ws.Range("F10").Value = 29.9

# Row 11 - OD10
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1

# Row 12 - OD20
$ws.Range("C12").Value = 1.1
$ws.Range("D12").Value = 1.2
$ws.Range("E12").Value = 1.4
$ws.Range("F12").Value = 1.5

# Row 13 - OD35
$ws.Range("C13").Value = 1.6
$ws.Range("D13").Value = 1.6
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1.9

# Row 14 - OD50
$ws.Range("C14").Value = 2.1
$ws.Range("D14").Value = 2.2
$ws.Range("E14").Value = 2.1
$ws.Range("F14").Value = 2.8

# Row 15 - OD85
$ws.Range("C15").Value = 3.2
$ws.Range("D15").Value = 3.2
$ws.Range("E15").Value = 3.1
$ws.Range("F15").Value = 3.8

# Apply the same header styling (bold + centered) used by the other header cells to the new ones
$ws.Range("E1:G1").Font.Bold = $true
$ws.Range("E1:G1").HorizontalAlignment = -4108

# Leave selection on G2, matching where the user ended up after adding the new data
$ws.Range("G2").Select()
